$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new entry "ChatVRM" to the next empty row in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($ws.Cells.Item(1,1).Value -eq $null) {
    $newRow = 1
} else {
    $newRow = $lastRow + 1
}

$ws.Cells.Item($newRow, 1).Value = "ChatVRM"

# Update selection to the newly added cell, matching the diff
$ws.Range("A" + $newRow).Select()
